# Auto-generated Excel COM-interop edit script
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per refreshed market data.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 3011.524
$ws.Range("I8").Value = 2171.2104
$ws.Range("K8").Value = 6513.6312
$ws.Range("M8").Value = -6374.6312
$ws.Range("H19").Value = 1683.7693
$ws.Range("I19").Value = 1400
$ws.Range("K19").Value = 1400
$ws.Range("M19").Value = -1225
$ws.Range("H62").Value = 4286.9062
$ws.Range("I62").Value = 3978
$ws.Range("J62").Value = 5213.625
$ws.Range("K62").Value = 3978
$ws.Range("L62").Value = 5213.625
$ws.Range("M62").Value = -3354
$ws.Range("N62").Value = -6461.625
$ws.Range("H65").Value = 4286.9062
$ws.Range("I65").Value = 3978
$ws.Range("J65").Value = 5213.625
$ws.Range("K65").Value = 19890
$ws.Range("L65").Value = 26068.125
$ws.Range("M65").Value = -16770
$ws.Range("N65").Value = -32308.125
$ws.Range("H97").Value = 1255.625
$ws.Range("J97").Value = 1255.625
$ws.Range("L97").Value = 3766.875
$ws.Range("N97").Value = -4758.875
$ws.Range("H98").Value = 1366.9286
$ws.Range("I98").Value = 1303.92
$ws.Range("K98").Value = 1303.92
$ws.Range("M98").Value = 194.0799999999999
$ws.Range("H111").Value = 7197093
$ws.Range("I111").Value = 10990134
$ws.Range("J111").Value = 152874.58
$ws.Range("K111").Value = 32970402
$ws.Range("L111").Value = 458623.74
$ws.Range("M111").Value = -32967335
$ws.Range("N111").Value = -464757.74
$ws.Range("H112").Value = 3166.9768
$ws.Range("I112").Value = 1422.25
$ws.Range("J112").Value = 3565.7715
$ws.Range("K112").Value = 4266.75
$ws.Range("L112").Value = 10697.3145
$ws.Range("M112").Value = -3158.75
$ws.Range("N112").Value = -12913.3145
$ws.Range("H122").Value = 1366.9286
$ws.Range("I122").Value = 1303.92
$ws.Range("K122").Value = 3911.76
$ws.Range("M122").Value = -1461.76
$ws.Range("H137").Value = 1928.1459
$ws.Range("I137").Value = 1629.3429
$ws.Range("J137").Value = 2732.6155
$ws.Range("K137").Value = 4888.028700000001
$ws.Range("L137").Value = 8197.8465
$ws.Range("M137").Value = -2338.028700000001
$ws.Range("N137").Value = -13297.8465
$ws.Range("H138").Value = 2296.79
$ws.Range("I138").Value = 896.58826
$ws.Range("J138").Value = 3018.106
$ws.Range("K138").Value = 2689.76478
$ws.Range("L138").Value = 9054.318000000001
$ws.Range("M138").Value = 2450.23522
$ws.Range("N138").Value = -19334.318

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14235.483
$ws.Range("I32").Value = 6773.8335
$ws.Range("J32").Value = 39818.285
$ws.Range("K32").Value = 6773.8335
$ws.Range("L32").Value = 39818.285
$ws.Range("M32").Value = -6486.8335
$ws.Range("N32").Value = -40392.285
$ws.Range("H45").Value = 84617680
$ws.Range("I45").Value = 110001740
$ws.Range("K45").Value = 110001740
$ws.Range("M45").Value = -110001363
$ws.Range("H61").Value = 9308.666999999999
$ws.Range("I61").Value = 9207.625
$ws.Range("K61").Value = 9207.625
$ws.Range("M61").Value = -8995.625
$ws.Range("H132").Value = 5586.9546
$ws.Range("I132").Value = 3347.9412
$ws.Range("J132").Value = 13199.6
$ws.Range("K132").Value = 10043.8236
$ws.Range("L132").Value = 39598.8
$ws.Range("M132").Value = -7513.8236
$ws.Range("N132").Value = -44658.8
$ws.Range("H136").Value = 9308.666999999999
$ws.Range("I136").Value = 9207.625
$ws.Range("K136").Value = 27622.875
$ws.Range("M136").Value = -25072.875

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("H20").Value = 2074.3044
$ws.Range("I20").Value = 2122.8125
$ws.Range("K20").Value = 2122.8125
$ws.Range("M20").Value = -1875.8125
$ws.Range("H86").Value = 4378.6333
$ws.Range("J86").Value = 5561.636
$ws.Range("L86").Value = 5561.636
$ws.Range("N86").Value = -7807.636
$ws.Range("H89").Value = 4378.6333
$ws.Range("J89").Value = 5561.636
$ws.Range("L89").Value = 27808.18
$ws.Range("N89").Value = -39040.18

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 1950
$ws.Range("I33").Value = 1950
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 1950
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -1571
$ws.Range("H94").Value = 1190
$ws.Range("I94").Value = 519.6
$ws.Range("J94").Value = 2028
$ws.Range("K94").Value = 519.6
$ws.Range("L94").Value = 2028
$ws.Range("M94").Value = -68.60000000000002
$ws.Range("N94").Value = -2930
$ws.Range("H122").Value = 1971.75
$ws.Range("I122").Value = 1995.1818
$ws.Range("J122").Value = 1714
$ws.Range("K122").Value = 5985.5454
$ws.Range("L122").Value = 5142
$ws.Range("M122").Value = -3535.5454
$ws.Range("N122").Value = -10042

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 4437.1816
$ws.Range("J75").Value = 4780.9
$ws.Range("L75").Value = 14342.7
$ws.Range("N75").Value = -16338.7
$ws.Range("H78").Value = 4437.1816
$ws.Range("J78").Value = 4780.9
$ws.Range("L78").Value = 43028.1
$ws.Range("N78").Value = -53012.1
$ws.Range("H103").Value = 662.375
$ws.Range("I103").Value = 424.75
$ws.Range("K103").Value = 1274.25
$ws.Range("M103").Value = -395.25
$ws.Range("H140").Value = 1983.7667
$ws.Range("J140").Value = 2167.36
$ws.Range("L140").Value = 6502.08
$ws.Range("N140").Value = -16862.08

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4857
$ws.Range("I80").Value = 3002.5
$ws.Range("K80").Value = 3002.5
$ws.Range("M80").Value = -2004.5
$ws.Range("H83").Value = 4857
$ws.Range("I83").Value = 3002.5
$ws.Range("K83").Value = 15012.5
$ws.Range("M83").Value = -10020.5
$ws.Range("H122").Value = 5456.4165
$ws.Range("I122").Value = 3298
$ws.Range("K122").Value = 9894
$ws.Range("M122").Value = -7444
$ws.Range("H126").Value = 2816.4773
$ws.Range("I126").Value = 2940.9722
$ws.Range("J126").Value = 2256.25
$ws.Range("K126").Value = 8822.9166
$ws.Range("L126").Value = 6768.75
$ws.Range("M126").Value = -6352.9166
$ws.Range("N126").Value = -11708.75

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1414.1428
$ws.Range("I16").Value = 1096
$ws.Range("K16").Value = 1096
$ws.Range("M16").Value = -926
$ws.Range("H61").Value = 2554.4736
$ws.Range("I61").Value = 1983.125
$ws.Range("J61").Value = 5601.6665
$ws.Range("K61").Value = 1983.125
$ws.Range("L61").Value = 5601.6665
$ws.Range("M61").Value = -1781.125
$ws.Range("N61").Value = -6005.6665
$ws.Range("H103").Value = 39999
$ws.Range("J103").Value = 39999
$ws.Range("L103").Value = 39999
$ws.Range("N103").Value = -42343
$ws.Range("H113").Value = 2554.4736
$ws.Range("I113").Value = 1983.125
$ws.Range("J113").Value = 5601.6665
$ws.Range("K113").Value = 1983.125
$ws.Range("L113").Value = 5601.6665
$ws.Range("M113").Value = 186.875
$ws.Range("N113").Value = -9941.666499999999

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 27000
$ws.Range("J74").Value = 29000
$ws.Range("L74").Value = 29000
$ws.Range("N74").Value = -30872
$ws.Range("H77").Value = 27000
$ws.Range("J77").Value = 29000
$ws.Range("L77").Value = 87000
$ws.Range("N77").Value = -96360
$ws.Range("H100").Value = 583.5454999999999
$ws.Range("I100").Value = 180
$ws.Range("K100").Value = 360
$ws.Range("M100").Value = 181
$ws.Range("H132").Value = 1930.2222
$ws.Range("I132").Value = 1633.9445
$ws.Range("J132").Value = 2522.7778
$ws.Range("K132").Value = 4901.833500000001
$ws.Range("L132").Value = 7568.3334
$ws.Range("M132").Value = -2371.833500000001
$ws.Range("N132").Value = -12628.3334
